$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Revenue growth multiplier assumption bumped from 1.34x to 1.44x (S6, and
# the shared formula spanning T6:V6 which repeats the same multiplier).
$ws.Range("S6").Formula = "=R6*1.44"
$ws.Range("T6:V6").Formula = "=S6*1.44"

# Leave the cell selection where the author left off (S16 on the frozen,
# bottom-right pane of the "model" sheet).
$ws.Activate()
$ws.Range("S16").Select()
